# Update symbol list (crypto price/volume snapshot refresh).
# Mirrors the upstream GitHub Actions job that re-scrapes coinranking.com
# and rewrites the Price (D) / Volume(1h) (E) columns for most rows, plus a
# coin-identity swap between rows 48 and 49 (Name/Link/Price/Volume).
#
# D/E store their content as literal text (t="inlineStr") in the workbook,
# not as numbers/percentages. A leading apostrophe forces Excel to keep the
# numeric-looking replacement text as plain text instead of reinterpreting it
# as a Number/Percentage. Excel then silently flags such a cell as "number
# stored as text" via an automatic quoted-text number format, so the style is
# reset back to "Normal" right after, one row at a time (contiguous D:E), to
# keep formatting identical to every other untouched cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'328.24"
$ws.Range("E2").Value = "'1.61%"
$ws.Range("D2:E2").Style = "Normal"

# Row 3 - OKB
$ws.Range("D3").Value = "'41.53"
$ws.Range("E3").Value = "'4.95%"
$ws.Range("D3:E3").Style = "Normal"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.612"
$ws.Range("E4").Value = "'-4.38%"
$ws.Range("D4:E4").Style = "Normal"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.08167"
$ws.Range("E5").Value = "'1.83%"
$ws.Range("D5:E5").Style = "Normal"

# Row 6 - FTXToken
$ws.Range("D6").Value = "'2.028"
$ws.Range("E6").Value = "'2.20%"
$ws.Range("D6:E6").Style = "Normal"

# Row 7 - KuCoinToken
$ws.Range("D7").Value = "'8.720"
$ws.Range("E7").Value = "'0.85%"
$ws.Range("D7:E7").Style = "Normal"

# Row 8 - GateToken
$ws.Range("D8").Value = "'4.521"
$ws.Range("E8").Value = "'-1.24%"
$ws.Range("D8:E8").Style = "Normal"

# Row 9 - BTSEToken
$ws.Range("D9").Value = "'2.942"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("D9:E9").Style = "Normal"

# Row 10 - MXToken
$ws.Range("D10").Value = "'0.9162"
$ws.Range("E10").Value = "'-1.35%"
$ws.Range("D10:E10").Style = "Normal"

# Row 11 - LiechtensteinCryptoassetsExchange
$ws.Range("D11").Value = "'0.1272"
$ws.Range("E11").Value = "'-0.64%"
$ws.Range("D11:E11").Style = "Normal"

# Row 12 - WazirX
$ws.Range("D12").Value = "'0.1954"
$ws.Range("E12").Value = "'-0.29%"
$ws.Range("D12:E12").Style = "Normal"

# Row 13 - MandalaExchangeToken
$ws.Range("D13").Value = "'0.09299"
$ws.Range("E13").Value = "'1.32%"
$ws.Range("D13:E13").Style = "Normal"

# Row 14 - BitrueCoin
$ws.Range("D14").Value = "'0.03767"
$ws.Range("E14").Value = "'5.77%"
$ws.Range("D14:E14").Style = "Normal"

# Row 15 - BitMartToken
$ws.Range("D15").Value = "'0.1059"
$ws.Range("E15").Value = "'1.22%"
$ws.Range("D15:E15").Style = "Normal"

# Row 16 - BitForexToken
$ws.Range("D16").Value = "'0.001301"
$ws.Range("E16").Value = "'0.82%"
$ws.Range("D16:E16").Style = "Normal"

# Row 17 - TigerCash
$ws.Range("D17").Value = "'0.006206"
$ws.Range("E17").Value = "'-1.03%"
$ws.Range("D17:E17").Style = "Normal"

# Row 19 - LEO
$ws.Range("D19").Value = "'3.440"
$ws.Range("E19").Value = "'2.70%"
$ws.Range("D19:E19").Style = "Normal"

# Row 21 - MCDex
$ws.Range("D21").Value = "'8.276"
$ws.Range("E21").Value = "'-4.97%"
$ws.Range("D21:E21").Style = "Normal"

# Row 22 - ProBitToken
$ws.Range("D22").Value = "'0.1394"
$ws.Range("E22").Value = "'1.64%"
$ws.Range("D22:E22").Style = "Normal"

# Row 23 - ZBToken
$ws.Range("E23").Value = "'-2.43%"
$ws.Range("E23").Style = "Normal"

# Row 24 - CoinExToken
$ws.Range("D24").Value = "'0.04431"
$ws.Range("E24").Value = "'0.58%"
$ws.Range("D24:E24").Style = "Normal"

# Row 25 - BitKan
$ws.Range("D25").Value = "'0.001259"
$ws.Range("E25").Value = "'-0.42%"
$ws.Range("D25:E25").Style = "Normal"

# Row 26 - HotbitToken
$ws.Range("D26").Value = "'0.004349"
$ws.Range("E26").Value = "'-1.02%"
$ws.Range("D26:E26").Style = "Normal"

# Row 27 - NitroEx
$ws.Range("D27").Value = "'0.0001181"
$ws.Range("E27").Value = "'3.53%"
$ws.Range("D27:E27").Style = "Normal"

# Row 39 - One
$ws.Range("D39").Value = "'0.02770"
$ws.Range("E39").Value = "'9.62%"
$ws.Range("D39:E39").Style = "Normal"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.05408"
$ws.Range("E40").Value = "'2.67%"
$ws.Range("D40:E40").Style = "Normal"

# Row 41 - KickToken
$ws.Range("D41").Value = "'0.007680"
$ws.Range("E41").Value = "'3.36%"
$ws.Range("D41:E41").Style = "Normal"

# Row 42 - BKEXToken
$ws.Range("D42").Value = "'0.1413"
$ws.Range("E42").Value = "'0.50%"
$ws.Range("D42:E42").Style = "Normal"

# Row 43 - Dexo
$ws.Range("D43").Value = "'0.008981"
$ws.Range("E43").Value = "'-6.60%"
$ws.Range("D43:E43").Style = "Normal"

# Row 44 - CEJI
$ws.Range("D44").Value = "'0.002123"
$ws.Range("E44").Value = "'0.20%"
$ws.Range("D44:E44").Style = "Normal"

# Row 45 - LocalTraders
$ws.Range("D45").Value = "'0.01157"
$ws.Range("E45").Value = "'15.89%"
$ws.Range("D45:E45").Style = "Normal"

# Row 46 - CoinLion
$ws.Range("D46").Value = "'0.00006356"
$ws.Range("E46").Value = "'-5.77%"
$ws.Range("D46:E46").Style = "Normal"

# Row 47 - Kangarootoken
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E47").Style = "Normal"

# Row 48 - BOLO -> CoinbaseStockToken
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002282"
$ws.Range("E48").Value = "'-0.41%"
$ws.Range("D48:E48").Style = "Normal"

# Row 49 - CoinbaseStockToken -> BOLO
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003252"
$ws.Range("E49").Value = "'8.30%"
$ws.Range("D49:E49").Style = "Normal"

# Row 50 - CryptobidCoin
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D50:E50").Style = "Normal"

# Row 51 - SpecialPowerGold
$ws.Range("E51").Value = "'0.07%"
$ws.Range("E51").Style = "Normal"
